$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the panel_query_time ("time_taken") column on the "data"
#        sheet -- this is a re-run of the same query a few hours later. ---
$timeTaken = @{
    2  = "2021-10-05 14:33:48.679138";
    3  = "2021-10-05 14:33:48.679146";
    4  = "2021-10-05 14:33:48.679149";
    5  = "2021-10-05 14:33:48.679152";
    6  = "2021-10-05 14:33:48.679155";
    7  = "2021-10-05 14:33:48.679157";
    8  = "2021-10-05 14:33:48.679160";
    9  = "2021-10-05 14:33:48.679163";
    10 = "2021-10-05 14:33:48.679165";
    11 = "2021-10-05 14:33:48.679168";
    12 = "2021-10-05 14:33:48.679171";
    13 = "2021-10-05 14:33:48.679173";
    14 = "2021-10-05 14:33:48.679176";
    15 = "2021-10-05 14:33:48.679179";
    16 = "2021-10-05 14:33:48.679181";
    17 = "2021-10-05 14:33:48.679184";
    18 = "2021-10-05 14:33:48.679187";
    19 = "2021-10-05 14:33:48.679189";
    20 = "2021-10-05 14:33:48.679192";
    21 = "2021-10-05 14:33:48.679195";
}
foreach ($row in $timeTaken.Keys) {
    $dataSheet.Cells.Item($row, 6).Value = $timeTaken[$row]
}

# --- 2. Add a new "metadata" tab right after "data", describing the
#        panelapp query that produced this workbook. ---

# Clone the header formatting (bold font + thin border + centered /
# top-aligned) from the existing "data" header row so the new sheet
# reuses the same cell style instead of registering a new one.
$dataSheet.Range("B1:F1").Copy()

$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Clone the index-column formatting (same style as B1:F1) onto A2.
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Ehlers Danlos syndromes"
$metaSheet.Range("C2").Value = 3180

# "1.0" must stay literal text (not become the number 1). Force it in
# as text via the quote-prefix, then re-paste the still-untouched A1's
# default formatting over it so the cell falls back to the default
# style (no quote-prefix / no explicit number format) while keeping
# the text value intact.
$metaSheet.Range("D2").Value = "'1.0"
$metaSheet.Range("A1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("E2").Value = "2020-10-28T02:00:33.256379Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:48.675599"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3180/?format=json"

# Restore the originally active sheet/selection.
$dataSheet.Activate()
